# Optimizasyon modellemesi matlab koduna başlandı. Parametreler tanıtıldı.
#
# The "p" parameter is moved from the Electromagnetic/Variable cell (C3)
# to the Electromagnetic/Dependent cell (E3):
#   C3: "α, em, w, p"                                   -> "α, em, w"
#   E3: "Tm, kw, Qs, Eph, Iph, Φpp, Nph, zQ, Dw, αs, Lm, Plossm"
#       -> "Tm, kw, Qs, p, Eph, Iph, Φpp, Nph, zQ, Dw, αs, Lm, Plossm"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E3 first so the pre-existing shared string is edited in place,
# then update C3 so its replacement string is the newly appended one -
# this mirrors the shared-string ordering produced by the real edit.
$ws.Range("E3").Value = "Tm, kw, Qs, p, Eph, Iph, Φpp, Nph, zQ, Dw, αs, Lm, Plossm"
$ws.Range("C3").Value = "α, em, w"

# Bring the view back to 100% zoom and move the active selection, matching
# the state the workbook was left in when it was re-saved.
$excel.ActiveWindow.Zoom = 100
$ws.Range("F11").Select()
